$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 541 - this shifts the existing rows 541-566
# down to 545-570 (and copies the column D date style down automatically,
# matching Excel's native row-insert formatting behaviour).
$ws.Rows("541:544").Insert()

# Common/fixed values shared by every data row in this table.
$mercadoId = 7
$mercado   = 'Terminal Hortofrutícola Agro Chillán'
$region    = 'Ñuble'
$codreg    = 16
$tipo      = 'Fruta'
$productoId = 100104
$producto  = 'Frutos de pepita'
$categoriaId = 100104002
$categoria = 'Manzana'
$unidad    = '$/caja 16 kilos empedrada'
$fecha     = 44509

$rows = @(
    @{ Row=541; K='Fuji royal';    L='Especial'; M=80;  N=11000; O=11000; P=11000; R='Provincia de Curicó'; S=688; T=16 },
    @{ Row=542; K='Fuji royal';    L='Primera';  M=160; N=9500;  O=10000; P=9750;  R='Provincia de Curicó'; S=609; T=16 },
    @{ Row=543; K='Granny Smith';  L='Primera';  M=120; N=9500;  O=10000; P=9750;  R='Provincia de Curicó'; S=609; T=16 },
    @{ Row=544; K='Granny Smith';  L='Segunda';  M=80;  N=8000;  O=8000;  P=8000;  R='Provincia de Curicó'; S=500; T=16 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $mercadoId
    $ws.Range("B$n").Value = $mercado
    $ws.Range("C$n").Value = $region
    $ws.Range("D$n").Value = $fecha
    $ws.Range("E$n").Value = $codreg
    $ws.Range("F$n").Value = $tipo
    $ws.Range("G$n").Value = $productoId
    $ws.Range("H$n").Value = $producto
    $ws.Range("I$n").Value = $categoriaId
    $ws.Range("J$n").Value = $categoria
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $unidad
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
}
